{"js": "// Update the worksheet date and each \"a\u00f7b=\" problem text to the new\n// values from the commit. Every source string is unique within the\n// document, so we resolve ALL search ranges against the pristine\n// (pre-edit) text first, then apply the text replacements \u2014 this avoids\n// a later search accidentally re-matching text that an earlier\n// replacement just wrote (e.g. \"39\u00f75=\" -> \"98\u00f72=\" followed by the\n// original \"98\u00f72=\" -> \"98\u00f77=\").\nconst pairs = [\n  [\"2025-11-25 Tuesday\", \"2025-11-26 Wednesday\"],\n  [\"46\u00f73=\", \"53\u00f74=\"],\n  [\"91\u00f73=\", \"92\u00f75=\"],\n  [\"39\u00f75=\", \"98\u00f72=\"],\n  [\"83\u00f79=\", \"97\u00f74=\"],\n  [\"48\u00f75=\", \"84\u00f72=\"],\n  [\"53\u00f73=\", \"18\u00f79=\"],\n  [\"90\u00f74=\", \"18\u00f75=\"],\n  [\"84\u00f74=\", \"15\u00f73=\"],\n  [\"29\u00f79=\", \"67\u00f72=\"],\n  [\"79\u00f78=\", \"72\u00f73=\"],\n  [\"81\u00f72=\", \"95\u00f77=\"],\n  [\"13\u00f74=\", \"24\u00f73=\"],\n  [\"59\u00f78=\", \"44\u00f72=\"],\n  [\"89\u00f78=\", \"49\u00f77=\"],\n  [\"58\u00f78=\", \"31\u00f77=\"],\n  [\"92\u00f74=\", \"84\u00f75=\"],\n  [\"24\u00f75=\", \"28\u00f76=\"],\n  [\"69\u00f76=\", \"76\u00f74=\"],\n  [\"84\u00f76=\", \"66\u00f74=\"],\n  [\"68\u00f73=\", \"93\u00f74=\"],\n  [\"98\u00f72=\", \"98\u00f77=\"],\n  [\"36\u00f75=\", \"88\u00f74=\"],\n  [\"63\u00f72=\", \"22\u00f74=\"],\n  [\"56\u00f79=\", \"71\u00f77=\"],\n  [\"22\u00f79=\", \"48\u00f79=\"],\n];\n\nconst body = context.document.body;\n\n// Kick off every search against the original text before anything is\n// mutated.\nconst searchResults = pairs.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// Now apply every replacement using the ranges captured above.\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText, newText] = pairs[i];\n  const items = searchResults[i].items;\n  if (items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${oldText}\", found ${items.length}`\n    );\n  }\n  items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and each \"a\u00f7b=\" division problem to the new\n# values from the commit.\n#\n# The document is a fixed-layout table: the date paragraph, then one\n# paragraph per filled table cell (blank \"spacer\" rows sit between each\n# row of problems). Walking $d.Paragraphs and writing to each target\n# paragraph's Range.Text by its (stable, 1-based) index updates every\n# value in place while preserving the run formatting, and \u2014 unlike a\n# sequential Find/Replace \u2014 it can't accidentally re-match text that an\n# earlier replacement just wrote (e.g. \"39\u00f75=\" -> \"98\u00f72=\" followed later\n# by the original \"98\u00f72=\" -> \"98\u00f77=\"), since each write targets a\n# specific paragraph object instead of searching document text.\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n$updates = @{\n    1   = \"2025-11-26 Wednesday\"\n    2   = \"53\u00f74=\"\n    3   = \"92\u00f75=\"\n    4   = \"98\u00f72=\"\n    5   = \"97\u00f74=\"\n    6   = \"84\u00f72=\"\n    26  = \"18\u00f79=\"\n    27  = \"18\u00f75=\"\n    28  = \"15\u00f73=\"\n    29  = \"67\u00f72=\"\n    30  = \"72\u00f73=\"\n    50  = \"95\u00f77=\"\n    51  = \"24\u00f73=\"\n    52  = \"44\u00f72=\"\n    53  = \"49\u00f77=\"\n    54  = \"31\u00f77=\"\n    74  = \"84\u00f75=\"\n    75  = \"28\u00f76=\"\n    76  = \"76\u00f74=\"\n    77  = \"66\u00f74=\"\n    78  = \"93\u00f74=\"\n    98  = \"98\u00f77=\"\n    99  = \"88\u00f74=\"\n    100 = \"22\u00f74=\"\n    101 = \"71\u00f77=\"\n    102 = \"48\u00f79=\"\n}\n\nforeach ($idx in $updates.Keys) {\n    $p = $paras.Item($idx)\n    $p.Range.Text = $updates[$idx]\n}\n"}
